# quarterly.xlsx update: roll the quarterly window forward one quarter.
# Drop the oldest quarter ("1399/06") from the header/labels and append the
# newest quarter ("1401/12") at the end (column N), shifting every other
# quarter's figures one column to the left (E<-F, F<-G, ... M<-N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels (row 8 "هزینه های عمومی و اداری" table, row 24 "تعداد پرسنل" table) ---
$quarters = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$headerRows = @(8, 24)
for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i   # column E = 5 .. N = 14
    foreach ($r in $headerRows) {
        $ws.Cells.Item($r, $col).Value = $quarters[$i]
    }
}

# --- Data rows: new (already-shifted) E:N figures for each line item ---
$dataRows = @{
    10 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)               # هزینه حمل و نقل و انتقال
    11 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)                # هزینه خدمات پس از فروش
    12 = @(0, 0, 0, 0, 0, 756449, 0, 132305, 58894, 266138)  # حق العمل و کمیسیون فروش
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)                # هزینه تبلیغات
    14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)                # هزینه مواد مصرفی
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)                # هزینه انرژی
    16 = @(2703, 2201, 3927, 4532, 5067, 5649, 6675, 7533, 7318, 8318)             # هزینه استهلاک
    17 = @(83467, 56931, 88297, 104467, 91761, 121497, 177706, 235284, 165646, 314759)  # هزینه حقوق و دستمزد
    18 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)                # هزینه مطالبات مشکوک الوصول
    19 = @(505827, 410291, 344045, 452162, 462066, -363663, 348031, 285930, 318637, 604629)  # سایر هزینه ها
    20 = @(591997, 469423, 436269, 561161, 558894, 519932, 532412, 661052, 550495, 1193844)  # جمع
    26 = @(272, 272, 415, 415, 415, 288, 288, 643, 2411, 743)          # تعداد پرسنل غیر تولیدی شرکت
    27 = @(1001, 1001, 890, 890, 890, 1081, 1081, 905, 911, 935)       # تعداد پرسنل تولیدی شرکت
}

foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i
        $ws.Cells.Item([int]$r, $col).Value = $vals[$i]
    }
}
